$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cell A1 = "Parameter", matching the bold/bordered header style used by B1:E1 ---
$ws.Range("A1").Value = "Parameter"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Bold the significant p-values (col E, p < 0.05) and the two fully-significant rows (10 & 12) ---
$ws.Range("E2").Font.Bold = $true
$ws.Range("E7").Font.Bold = $true
$ws.Range("E8").Font.Bold = $true
$ws.Range("B10:E10").Font.Bold = $true
$ws.Range("E11").Font.Bold = $true
$ws.Range("B12:E12").Font.Bold = $true
$ws.Range("E13").Font.Bold = $true
$ws.Range("E14").Font.Bold = $true
$ws.Range("E15").Font.Bold = $true
$ws.Range("E16").Font.Bold = $true

# --- 3. Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- 4. Final selection, as last left by the author ---
$ws.Range("C21").Select() | Out-Null

Write-Output "done"
